# feat (paper info): update name data
#
# Row 28 (paper id 26) lists the author of "University of Moratuwa, Sri
# Lanka / Al Baha University, Saudi Arabia" under the spelled-out name
# "Mohamed Fazil Mohamed Firdhous". The accepted-papers sheet is updated
# to use the author's preferred short form "M F M Firdhous" (second
# co-author "Rahmat Budiarto" is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("D28").Value = "M F M Firdhous`nRahmat Budiarto"

# Reflect the author's final on-sheet cursor position/selection from the
# edit session.
$ws.Range("G26").Select()
